$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update query timestamps on the "data" sheet (column F) ---
$dataSheet.Cells.Item(2,6).Value = "2021-10-05 14:22:21.384979"
$dataSheet.Cells.Item(3,6).Value = "2021-10-05 14:22:21.384999"
$dataSheet.Cells.Item(4,6).Value = "2021-10-05 14:22:21.385002"
$dataSheet.Cells.Item(5,6).Value = "2021-10-05 14:22:21.385004"
$dataSheet.Cells.Item(6,6).Value = "2021-10-05 14:22:21.385026"

# --- Add a new "metadata" worksheet right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (row 1), columns B..G
$ws.Cells.Item(1,2).Value = "data_name"
$ws.Cells.Item(1,3).Value = "data_id"
$ws.Cells.Item(1,4).Value = "data_version"
$ws.Cells.Item(1,5).Value = "data_version_created"
$ws.Cells.Item(1,6).Value = "panel_query_time"
$ws.Cells.Item(1,7).Value = "panel_get_request"

# Copy the bold/bordered header style from the "data" sheet header row
$dataSheet.Range("B1:F1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row (row 2)
$ws.Cells.Item(2,1).Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Cells.Item(2,2).Value = "Prostate cancer pertinent cancer susceptibility"
$ws.Cells.Item(2,3).Value = 17

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "1.1"

$ws.Cells.Item(2,5).Value = "2019-06-20T15:13:53.649658Z"
$ws.Cells.Item(2,6).Value = "2021-10-05 14:22:21.382612"
$ws.Cells.Item(2,7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/17/?format=json"

$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null
